# Mise a jour du fichier via Shiny
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("pro")
$updates = @(
    @(2, 990321.51237977296),
    @(3, 1070415.4804183352),
    @(4, 728301.5619912335),
    @(5, 834178.08056133857),
    @(6, 1087372.5359733377),
    @(7, 762639.08575931063),
    @(8, 670980.12225690088),
    @(9, 676210.77560273965),
    @(10, 1073428.1063382323),
    @(11, 1170549.5376254548),
    @(12, 1404638.3305841736),
    @(13, 1260018.8242749723),
    @(14, 1326740.5622562943),
    @(15, 1221894.191260132),
    @(16, 1070946.9989942736),
    @(17, 942444.29300389811),
    @(18, 1063111.2264746665),
    @(19, 1151682.7509259763),
    @(20, 1101630.9827816947),
    @(21, 1275255),
    @(22, 1098305),
    @(23, 1098305),
    @(24, 1237678.4097270947),
    @(25, 1473204.2116241555),
    @(26, 1095229.303664732)
)
foreach ($u in $updates) {
    $ws.Cells.Item($u[0], 2).Value = $u[1]
}

$ws = $wb.Worksheets.Item("ind")
$updates = @(
    @(2, 1000107.9307687287),
    @(3, 813274.58106468047),
    @(4, 626441.23136063211),
    @(5, 734511.50226787559),
    @(6, 758323.59585760732),
    @(7, 840750.07366821694),
    @(8, 923176.55147882644),
    @(9, 917681.45295811887),
    @(10, 486470.23831790633),
    @(11, 632751.49879112281),
    @(12, 676976.06591093249),
    @(13, 544302.36455150344),
    @(14, 972632.43816519203),
    @(15, 939661.84704094846),
    @(16, 906691.25591670431),
    @(17, 879215.76331316784),
    @(18, 1246685.9071956915),
    @(19, 1390953.0536615688),
    @(20, 1050700.3497326132),
    @(21, 1132360.9986755624),
    @(22, 631936.32988133945),
    @(23, 679560.51706080278),
    @(24, 771145.49240592425),
    @(25, 879215.76331316796),
    @(26, 756313.3139223305),
    @(27, 681115.81425892492),
    @(28, 685454.13154719851),
    @(29, 482999.32476110605),
    @(30, 615451.03431921755),
    @(31, 657580.12297797366),
    @(32, 664906.92100558325),
    @(33, 978127.53668589925),
    @(34, 952483.74358926539),
    @(35, 923176.55147882609),
    @(36, 892037.65986148512),
    @(37, 1089861.2066069476),
    @(38, 870057.26577865588),
    @(39, 1245555.6646936545),
    @(40, 1077039.3100586308),
    @(41, 1056890.615482704),
    @(42, 1178515.4627410255),
    @(43, 1182728.3716069008),
    @(44, 1005969.3691908162),
    @(45, 1138767.5834412426),
    @(46, 920063.72033100191),
    @(47, 1072679.8052997568),
    @(48, 1125291.8929049054),
    @(49, 884627.93483109295),
    @(50, 699749.59555758315),
    @(51, 1144815.3668506786),
    @(52, 1105801.8391470877),
    @(53, 1264249.4610423064),
    @(54, 968445.68615854834),
    @(55, 823633.73591681162),
    @(56, 1176426.1867035949),
    @(57, 945930.92577339546),
    @(58, 795098.61670373112),
    @(59, 916653.20835892786),
    @(60, 907382.85334839148),
    @(61, 963249.70343915641),
    @(62, 431427.54037674918),
    @(63, 369966.53573376074),
    @(64, 751449.32714093931),
    @(65, 746289.71330065699),
    @(66, 720871.09224771557),
    @(67, 1081164.9894799686),
    @(68, 922290.22699134727),
    @(69, 953250.14758077799),
    @(70, 683534.08508741891),
    @(71, 870506.00799628056),
    @(72, 1230548.1593332535),
    @(73, 1199379.9951734564),
    @(74, 681091.98031742696),
    @(75, 841382.46239189501),
    @(76, 1017309.0362030774),
    @(77, 926087.69392797025),
    @(78, 1101113.6546200849),
    @(79, 977314.11923495145),
    @(80, 1095856.0032444424),
    @(81, 825716.22290052124),
    @(82, 797890.91839350585),
    @(83, 1100285.8594182963),
    @(84, 1014009.6871968094),
    @(85, 929522.75785041205),
    @(86, 855870.07286193897),
    @(87, 988135.92875113327),
    @(88, 1002867.0304461966),
    @(89, 994836.19079975481),
    @(90, 1063950.9477360607),
    @(91, 1069774.433664954),
    @(92, 1117483.4168141712),
    @(93, 1071218.9649397545),
    @(94, 1220653.3189890692),
    @(95, 1333519.5802183093),
    @(96, 1315137.49670034),
    @(97, 1325900.1869863209),
    @(98, 1301124.2243798147),
    @(99, 1176888.9193337392),
    @(100, 1215265.4916465275),
    @(101, 1203895.2144892709)
)
foreach ($u in $updates) {
    $ws.Cells.Item($u[0], 2).Value = $u[1]
}

$ws = $wb.Worksheets.Item("conso")
$updates = @(
    @(2, 893891.51093727164),
    @(3, 966185.40929496987),
    @(4, 657383.77010683052),
    @(5, 752951.3210933503),
    @(6, 981491.54481152759),
    @(7, 688379.03980535059),
    @(8, 605645.68933115818),
    @(9, 610366.91763670254),
    @(10, 968908.26941384177),
    @(11, 1056573.6405973702),
    @(12, 1267868.4013958124),
    @(13, 1137331.2884043553),
    @(14, 1197555.5984802905),
    @(15, 1102918.226054572),
    @(16, 966668.9240178609),
    @(17, 850678.54938435135),
    @(18, 959596.1616547599),
    @(19, 1039543.9330991084),
    @(20, 994365.63566010527),
    @(21, 1151084),
    @(22, 1198807),
    @(23, 1198807),
    @(24, 724325.95421133691),
    @(25, 862162.60859564901),
    @(26, 640960.53079902742)
)
foreach ($u in $updates) {
    $ws.Cells.Item($u[0], 2).Value = $u[1]
}

# Restore the final selection state recorded for each sheet (D87, with "pro" left as the active tab)
foreach ($name in @("ind","VA","conso","pro")) {
    $sel = $wb.Worksheets.Item($name)
    $sel.Range("D87").Select()
}
